$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = "SOC_XSIT_5"
$ws.Cells.Item(6, 2).Value = 41851
$ws.Cells.Item(6, 3).Value = 39957
$ws.Cells.Item(6, 5).Value = "F"
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(6, 7).Value = "social"
$ws.Cells.Item(6, 8).Value = 5

$ws.Cells.Item(7, 1).Value = "SOC_XSIT_6"
$ws.Cells.Item(7, 2).Value = 41851
$ws.Cells.Item(7, 3).Value = 39968
$ws.Cells.Item(7, 5).Value = "M"
$ws.Cells.Item(7, 6).Value = 5
$ws.Cells.Item(7, 8).Value = 5

$ws.Cells.Item(8, 1).Value = "SOC_XSIT_7"
$ws.Cells.Item(8, 2).Value = 41851
$ws.Cells.Item(8, 3).Value = 40147
$ws.Cells.Item(8, 5).Value = "F"
$ws.Cells.Item(8, 6).Value = 4.6
$ws.Cells.Item(8, 8).Value = 4

$ws.Cells.Item(9, 1).Value = "SOC_XSIT_8"
$ws.Cells.Item(9, 2).Value = 41851
$ws.Cells.Item(9, 3).Value = 40372
$ws.Cells.Item(9, 5).Value = "F"
$ws.Cells.Item(9, 6).Value = 3.11
$ws.Cells.Item(9, 8).Value = 3

$ws.Cells.Item(10, 1).Value = "SOC_XSIT_9"
$ws.Cells.Item(10, 2).Value = 41851
$ws.Cells.Item(10, 3).Value = 40315
$ws.Cells.Item(10, 5).Value = "M"
$ws.Cells.Item(10, 6).Value = 4.1
$ws.Cells.Item(10, 8).Value = 4

$ws.Cells.Item(11, 1).Value = "SOC_XSIT_10"
$ws.Cells.Item(11, 2).Value = 41851
$ws.Cells.Item(11, 3).Value = 40386
$ws.Cells.Item(11, 5).Value = "M"
$ws.Cells.Item(11, 6).Value = 4
$ws.Cells.Item(11, 8).Value = 4

$ws.Cells.Item(12, 1).Value = "SOC_XSIT_11"
$ws.Cells.Item(12, 2).Value = 41851
$ws.Cells.Item(12, 3).Value = 40289
$ws.Cells.Item(12, 5).Value = "F"
$ws.Cells.Item(12, 6).Value = 4.2
$ws.Cells.Item(12, 8).Value = 4

$ws.Cells.Item(13, 1).Value = "SOC_XSIT_12"
$ws.Cells.Item(13, 2).Value = 41852
$ws.Cells.Item(13, 3).Value = 40339
$ws.Cells.Item(13, 5).Value = "F"
$ws.Cells.Item(13, 6).Value = 4
$ws.Cells.Item(13, 8).Value = 4

$ws.Cells.Item(14, 1).Value = "SOC_XSIT_13"
$ws.Cells.Item(14, 2).Value = 41852
$ws.Cells.Item(14, 3).Value = 40304
$ws.Cells.Item(14, 5).Value = "M"
$ws.Cells.Item(14, 6).Value = 4.1
$ws.Cells.Item(14, 8).Value = 4

$ws.Cells.Item(15, 1).Value = "SOC_XSIT_14"
$ws.Cells.Item(15, 2).Value = 41852
$ws.Cells.Item(15, 3).Value = 39995
$ws.Cells.Item(15, 5).Value = "F"
$ws.Cells.Item(15, 6).Value = 4.11
$ws.Cells.Item(15, 8).Value = 4

$ws.Cells.Item(16, 1).Value = "SOC_XSIT_15"
$ws.Cells.Item(16, 2).Value = 41852
$ws.Cells.Item(16, 3).Value = 39948
$ws.Cells.Item(16, 5).Value = "F"
$ws.Cells.Item(16, 6).Value = 5.1
$ws.Cells.Item(16, 8).Value = 5

$ws.Cells.Item(17, 1).Value = "SOC_XSIT_16"
$ws.Cells.Item(17, 2).Value = 41852
$ws.Cells.Item(17, 3).Value = 39971
$ws.Cells.Item(17, 5).Value = "M"
$ws.Cells.Item(17, 6).Value = 5
$ws.Cells.Item(17, 8).Value = 5

$ws.Cells.Item(18, 1).Value = "SOC_XSIT_17"
$ws.Cells.Item(18, 2).Value = 41855
$ws.Cells.Item(18, 3).Value = 40421
$ws.Cells.Item(18, 5).Value = "M"
$ws.Cells.Item(18, 6).Value = 3.9
$ws.Cells.Item(18, 8).Value = 3

$ws.Range("A7").Select()
